# Update workbook to reflect data through 2022-05-23 (commit: "Add data for 2022-05-31")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (title tracks the "through" date)
$ws.Name = "Through 2022-05-23"

# Update the label for the May row to reflect the new "through" date
$ws.Range("A6").Value = "May (through 05-23)"

# Update May row (row 6) values
$ws.Range("C6").Value = 36
$ws.Range("D6").Value = 46
$ws.Range("E6").Value = 33
$ws.Range("F6").Value = 34
$ws.Range("G6").Value = 45
$ws.Range("H6").Value = 87
$ws.Range("I6").Value = 83

# Update Total row (row 7) values
$ws.Range("C7").Value = 198
$ws.Range("D7").Value = 299
$ws.Range("E7").Value = 279
$ws.Range("F7").Value = 189
$ws.Range("G7").Value = 307
$ws.Range("H7").Value = 610
$ws.Range("I7").Value = 635
